$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Try copying format from an existing similar block (rows 83-85, header+table) to rows 89-91
$ws.Range("A83:C83").Copy()
$ws.Range("A89:C89").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A89").Value = "// Here is how we count across columns"
$ws.Range("A89:C89").Merge()

$ws.Range("A60:D60").Copy()
$ws.Range("A90:D90").PasteSpecial(-4122)
$ws.Range("A90").Value = "COUNTIF"
$ws.Range("B90").Value = "XX_COL, Value"
$ws.Range("C90").Value = "XX_COL, Value2"
$ws.Range("D90").Value = "Here, each record's value at column XX is compared to ""Value"" . If the values match, then a counter is incremented. You can have as many columns as you need"

$ws.Range("A62:D62").Copy()
$ws.Range("A91:D91").PasteSpecial(-4122)
$ws.Range("D91").Value = "Only 2 columns shown here, but can have as many as neeeded."
